$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 corresponds to the "12 kOhm Resistor (2.2V)" component.
# Update the Mouser part number, unit price, and description to the new part.
$ws.Range("C13").Value = "667-ERA-3AEB123V"
$ws.Range("E13").Value = 0.299
$ws.Range("F13").Value = "Thin Film Resistors - SMD 0603 1/10W 12Kohms"

# Reflect the new selection left behind after the edit.
$ws.Range("E27").Select()
